$p = $ppt.ActivePresentation

# Slide 5: "Arquitectura & Tecnologias" - TextBox 3 (shape index 3)
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(3)
$sh5.TextFrame.TextRange.Text = "🏗️ TIPO DE ARQUITECTURA: MICROSERVICIOS CON API GATEWAY`r  • API Gateway centralizado (Express.js en Node.js)`r  • 2 servicios backend independientes (desacoplados)`r  • 2 bases de datos especializadas (poliglot persistence)`r  • Frontend desvinculado de la lógica del servidor`r  • Cada servicio puede escalarse independientemente`r`r🔧 STACK TECNOLÓGICO:`r  API Gateway: Express.js en Node.js (puerto 8080)`r  Backend 1: FastAPI + Motor (async MongoDB) en Python`r  Backend 2: Express.js + MySQL en Node.js`r  Frontend: React 19 + Vite + React Router v7`r  DevOps: Docker + Docker Compose + Nginx`r`r🌐 API GATEWAY - PUNTO DE ENTRADA ÚNICO:`r  • Autenticación JWT centralizada`r  • Proxy inverso inteligente hacia servicios backend`r  • Mapeo de rutas: /auth/* → /users/* (users service)`r  • Manejo de rutas públicas y protegidas`r  • Logging centralizado con Morgan`r  • CORS configurado para frontend`r`r📱 ARQUITECTURA DE FRONTEND:`r  • SPA (Single Page Application) con React`r  • CSR (Client-Side Rendering)`r  • Enrutamiento con React Router v7`r  • Comunicación con gateway en puerto 8080`r`r🔌 FLUJO DE COMUNICACIÓN:`r  Frontend (5173) → API Gateway (8080) → Servicios (8000, 3001)"

# Slide 6: "Lecciones & Aspectos Avanzados" - TextBox 3 (shape index 3)
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(3)
$sh6.TextFrame.TextRange.Text = "🎓 LECCIONES APRENDIDAS:`r  • Integración eficiente de múltiples lenguajes (Python + Node.js)`r  • Orquestación compleja con Docker Compose (6 contenedores)`r  • Gestión de dependencias entre servicios`r  • Implementación de API Gateway para centralizar autenticación`r`r⚡ ASPECTOS AVANZADOS IMPLEMENTADOS:`r  • API Gateway con proxy inverso inteligente (express-http-proxy)`r  • Mapeo de rutas y validación centralizada de JWT`r  • Async/await en FastAPI: operaciones no-bloqueantes`r  • Motor: driver async para MongoDB con alta concurrencia`r  • JWT con tokens con expiración y revocación`r  • CORS configurado en gateway y backends para seguridad`r  • Rutas protegidas: ProtectedRoute y AdminRoute en React`r  • Validación de entrada con Pydantic schemas`r  • Hashing seguro de contraseñas con bcrypt`r  • Control de acceso basado en roles (RBAC)`r`r⭐ POR QUÉ MERECE BUENA NOTA:`r  ✅ Sistema COMPLETO y FUNCIONAL (Full-Stack + Gateway)`r  ✅ Arquitectura PROFESIONAL (3-tier con API Gateway)`r  ✅ SEGURIDAD en múltiples niveles (JWT + RBAC + CORS)`r  ✅ BUENAS PRÁCTICAS y patrones de diseño actuales`r  ✅ ESCALABLE y MANTENIBLE para producción`r  ✅ Containerizado con Docker para reproducibilidad"
